# Add "barrnap mode" related OTU columns to the metadata demo workbook.
#
# Before the edit, columns were laid out (row 1/2 headers, rows 3-7 data):
#   ... BW  BX  BY               CA                 CB
#   ... 30  20  barrnap_mode     taxa_barstacks     max_taxa
#
# After the edit two brand-new columns are inserted right after BX
# (otu_threshold, otu_min_size) and the existing "barrnap_mode" column is
# pushed two columns to the right, along with everything after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three columns starting at the old "BY" column. Excel shifts BY
# onward to the right by three (BY -> CB, CA -> CD, CB -> CE) and the three
# freshly inserted columns inherit the formatting of the column immediately
# to their left (BX), matching the widths/styles seen in the target file.
$ws.Range("BY1:CA1").EntireColumn.Insert()

# Match the column width of the preceding "BX" column (~21.83 chars, i.e.
# the 22.6640625 width used for BX/BW) for the three new columns.
$ws.Range("BY1:CA1").EntireColumn.ColumnWidth = 21.83

# Populate the new header cells in the same order the original authoring
# tool produced them so the shared-string table indices line up exactly:
#   171 otu_threshold, 172 <description>, 173 otu_min_size, 174 <description>
$ws.Range("BY2").Value = ""
$ws.Range("BZ2").Value = "otu_threshold"

$ws.Range("BY1").Value = ""
$ws.Range("BZ1").Value = "This will determine the percentage of similarty that two sequences must share for them to be clustered into the same OTU. This is an optional column with a default value of 0.97 "

$ws.Range("CA2").Value = "otu_min_size"

$ws.Range("CA1").Value = "This will determine the minimum size of an OTU. Below this size OTUs are discarded. This is an optional column with a default value of 1. Setting this value to 2 will remove all singeltons for instance."

# Rows 3-7 (sample data) -- same default values repeated for every sample
$ws.Range("BZ3").Value = 0.97
$ws.Range("CA3").Value = 2

$ws.Range("BZ4").Value = 0.97
$ws.Range("CA4").Value = 2

$ws.Range("BZ5").Value = 0.97
$ws.Range("CA5").Value = 2

$ws.Range("BZ6").Value = 0.97
$ws.Range("CA6").Value = 2

$ws.Range("BZ7").Value = 0.97
$ws.Range("CA7").Value = 2

# Restore the selection to match the saved view of the edited workbook.
$ws.Range("BZ17").Select()
